# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps on the
# per-language report sheets (zh-cn and de-de) to reflect the
# newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-08 06:10:41"
$wsZhCn.Range("G4").Value = "2016-03-08 06:10:58"
$wsZhCn.Range("D5").Value = "2016-03-08 06:10:41"
$wsZhCn.Range("G5").Value = "2016-03-08 06:10:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-08 06:10:45"
$wsDeDe.Range("G4").Value = "2016-03-08 06:11:03"
$wsDeDe.Range("D5").Value = "2016-03-08 06:10:45"
$wsDeDe.Range("G5").Value = "2016-03-08 06:11:03"
